$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "LF"
$ws.Range("A3").Value = "LF Lag"
$ws.Range("B2").Value = "-0.326***"
$ws.Range("C2").Value = "'0.159"
$ws.Range("B3").Value = "-0.261*"
$ws.Range("C3").Value = "-0.454**"
